$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reposition the workbook window (xWindow 1900 -> 12290) to match the saved view state.
$excel.ActiveWindow.Left = 12290

# Update the Pre_ISI (column I) values for the rows whose timing was normalised to 3.01.
$rows = @(2, 7, 11, 15, 19, 21, 30, 42, 45, 50, 53, 57, 71, 84, 96, 101, 106, 107)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = 3.01
}

# Flip the High_Prob_Choice_Corr flag (column G, row 4) from 0 to 1.
$ws.Cells.Item(4, 7).Value = 1

# Move the sheet's active cell/selection from Q9 to M8.
[void]$ws.Range("M8").Select()
